$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MML_sig_config")

$ws.Range("A1").Value = 'username'
$ws.Range("B1").Value = 'aspect'
$ws.Range("C1").Value = 'panel_switch'
$ws.Range("D1").Value = 'properties'

$ws.Range("A2").Value = 'SG 1R'
$ws.Range("B2").Value = '4-distant-approach'
$ws.Range("D2").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;SigType:DistantApproach}'

$ws.Range("A3").Value = 'SG 2R'
$ws.Range("B3").Value = '3-home-diverging'
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;SigType:Home;Diverging:Right}'

$ws.Range("A4").Value = 'SG 3R'
$ws.Range("B4").Value = '2-mainline-starter'
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML}'

$ws.Range("A5").Value = 'SG 4R'
$ws.Range("B5").Value = '2-mainline-starter'
$ws.Range("D5").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;SigType:LSS}'

$ws.Range("A6").Value = 'SG 5L'
$ws.Range("B6").Value = '2-mainline-starter'
$ws.Range("D6").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML;SigType:LSS}'

$ws.Range("A7").Value = 'SG 6L'
$ws.Range("B7").Value = '2-mainline-starter'
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML}'

$ws.Range("A8").Value = 'SG 7L'
$ws.Range("B8").Value = '3-home-diverging'
$ws.Range("C8").Value = $true
$ws.Range("D8").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML;SigType:Home;Diverging:Left}'

$ws.Range("A9").Value = 'SG 8L'
$ws.Range("B9").Value = '4-distant-approach'
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML;SigType:DistantApproach}'

$ws.Range("A10").Value = 'SG 9R'
$ws.Range("B10").Value = '3-home'
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;SigType:Home}'

$ws.Range("A11").Value = 'SG 22R'
$ws.Range("B11").Value = '4-distant-approach'
$ws.Range("D11").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;SigType:DistantApproach}'

$ws.Range("A12").Value = 'SG 10R'
$ws.Range("B12").Value = '3-home-diverging'
$ws.Range("C12").Value = $true
$ws.Range("D12").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;SigType:Home;Diverging:Right}'

$ws.Range("A13").Value = 'SG 11R'
$ws.Range("B13").Value = '2-mainline-starter'
$ws.Range("C13").Value = $true
$ws.Range("D13").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;SigType:Starter}'

$ws.Range("A14").Value = 'SG 12R'
$ws.Range("B14").Value = '3-home-diverging'
$ws.Range("C14").Value = $true
$ws.Range("D14").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;Diverging:Left}'

$ws.Range("A15").Value = 'SG 13R'
$ws.Range("B15").Value = '2-mainline-starter'
$ws.Range("D15").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;SigType:LSS}'

$ws.Range("A16").Value = 'SG 14L'
$ws.Range("B16").Value = '2-mainline-starter'
$ws.Range("C16").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("D16").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML;SigType:LSS}'

$ws.Range("A17").Value = 'SG 15L'
$ws.Range("B17").Value = '3-home-diverging'
$ws.Range("C17").Value = $true
$ws.Range("D17").Value = '{OnTrack:MNL;Direction:Right;StationCode:MML;Diverging:Right}'

$ws.Range("A18").Value = 'SG 16L'
$ws.Range("B18").Value = '2-mainline-starter'
$ws.Range("C18").Value = $true
$ws.Range("D18").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML;SigType:Starter}'

$ws.Range("A19").Value = 'SG 17L'
$ws.Range("B19").Value = '3-home-diverging'
$ws.Range("C19").Value = $true
$ws.Range("D19").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML;Diverging:Left}'

$ws.Range("A20").Value = 'SG 18L'
$ws.Range("B20").Value = '3-home'
$ws.Range("C20").Value = $true
$ws.Range("D20").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML;SigType:Home}'

$ws.Range("A21").Value = 'SG 19L'
$ws.Range("B21").Value = '4-distant-approach'
$ws.Range("C21").Value = $true
$ws.Range("D21").Value = '{OnTrack:MNL;Direction:Left;StationCode:MML;SigType:DistantApproach}'

$ws.Range("A22").Value = 'SG 20R'
$ws.Range("B22").Value = '2-general'
$ws.Range("C22").Value = $true
$ws.Range("D22").Value = '{OnTrack:LPL;Direction:Right;StationCode:MML}'

$ws.Range("A23").Value = 'SG 21L'
$ws.Range("B23").Value = '2-general'
$ws.Range("C23").Value = $true
$ws.Range("D23").Value = '{OnTrack:LPL;Direction:Left;StationCode:MML}'

$ws.Columns.Item(2).ColumnWidth = 16.72
$ws.Columns.Item(3).ColumnWidth = 12.2

$ws.Range("A24").Select()

$ws.Activate()

